# Edit matching the commit diff: shift the MSI/MSIDN/SERIAL sample rows down by
# one "generation" (decrement the trailing digit of each serial-like value)
# and drop the last data row (row 13), leaving two data rows instead of three.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the remaining two data rows with their new values (column-major
# order so new shared-string entries land in the same order Excel would
# have produced while typing column-by-column).
$ws.Range("B11").Value = "8957732111198172291"
$ws.Range("B12").Value = "8957732111198172290"

$ws.Range("C11").Value = "3016875982"
$ws.Range("C12").Value = "3016875893"

$ws.Range("D11").Value = "732111198172291"
$ws.Range("D12").Value = "732111198172290"

# Remove the old third data row entirely (was row 13).
$ws.Rows.Item(13).Delete() | Out-Null
